# Update Name of Algo
# Applies new KNN-imputed values to specific cells, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 6.598999999999999
$ws.Range("D5").Value  = -8.148999999999999
$ws.Range("E7").Value  = 13.078
$ws.Range("D9").Value  = -7.634
$ws.Range("D11").Value = -8.171000000000001
$ws.Range("E11").Value = 12.949
$ws.Range("B21").Value = 6.434
$ws.Range("D21").Value = -7.7
$ws.Range("E21").Value = 12.054
$ws.Range("B23").Value = 6.945
$ws.Range("B25").Value = 6.556999999999999
